$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(155, 2).Value = 6810132
$ws.Cells.Item(155, 6).Value = 'SintTruidense'
$ws.Cells.Item(155, 7).Value = 'Charleroi'
$ws.Cells.Item(155, 8).Value = 1
$ws.Cells.Item(155, 9).Value = 0
$ws.Cells.Item(155, 10).Value = 'H'
$ws.Cells.Item(155, 11).Value = 2.2
$ws.Cells.Item(155, 12).Value = 3.4
$ws.Cells.Item(155, 13).Value = 3.2
$ws.Cells.Item(155, 14).Value = 2.3
$ws.Cells.Item(155, 15).Value = 3.3
$ws.Cells.Item(155, 16).Value = 3
$ws.Cells.Item(155, 17).Value = -0.25
$ws.Cells.Item(155, 18).Value = 2
$ws.Cells.Item(155, 19).Value = 1.85
$ws.Cells.Item(155, 20).Value = 2.25
$ws.Cells.Item(155, 21).Value = 1.85
$ws.Cells.Item(155, 22).Value = 2
$ws.Cells.Item(155, 23).Value = 1.3
$ws.Cells.Item(155, 24).Value = -1
$ws.Cells.Item(155, 26).Value = 1
$ws.Cells.Item(155, 27).Value = -1
$ws.Cells.Item(155, 28).Value = -1
$ws.Cells.Item(155, 29).Value = 1
$ws.Cells.Item(156, 2).Value = 6810130
$ws.Cells.Item(156, 6).Value = 'Antwerp'
$ws.Cells.Item(156, 7).Value = 'Westerlo'
$ws.Cells.Item(156, 8).Value = 2
$ws.Cells.Item(156, 9).Value = 2
$ws.Cells.Item(156, 10).Value = 'D'
$ws.Cells.Item(156, 11).Value = 1.363
$ws.Cells.Item(156, 12).Value = 5
$ws.Cells.Item(156, 13).Value = 7.5
$ws.Cells.Item(156, 14).Value = 1.333
$ws.Cells.Item(156, 15).Value = 5.25
$ws.Cells.Item(156, 16).Value = 8
$ws.Cells.Item(156, 17).Value = -1.5
$ws.Cells.Item(156, 18).Value = 1.9
$ws.Cells.Item(156, 19).Value = 1.95
$ws.Cells.Item(156, 20).Value = 3
$ws.Cells.Item(156, 21).Value = 1.8
$ws.Cells.Item(156, 22).Value = 2.05
$ws.Cells.Item(156, 23).Value = -1
$ws.Cells.Item(156, 24).Value = 4.25
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 27).Value = 0.95
$ws.Cells.Item(156, 28).Value = 0.8
$ws.Cells.Item(156, 29).Value = -1
$ws.Cells.Item(159, 2).Value = 6810142
$ws.Cells.Item(159, 6).Value = 'OH Leuven'
$ws.Cells.Item(159, 7).Value = 'Eupen'
$ws.Cells.Item(159, 8).Value = 3
$ws.Cells.Item(159, 9).Value = 0
$ws.Cells.Item(159, 10).Value = 'H'
$ws.Cells.Item(159, 11).Value = 1.75
$ws.Cells.Item(159, 12).Value = 4
$ws.Cells.Item(159, 13).Value = 4
$ws.Cells.Item(159, 14).Value = 1.8
$ws.Cells.Item(159, 15).Value = 3.8
$ws.Cells.Item(159, 16).Value = 3.8
$ws.Cells.Item(159, 17).Value = -0.5
$ws.Cells.Item(159, 18).Value = 1.825
$ws.Cells.Item(159, 19).Value = 2.025
$ws.Cells.Item(159, 20).Value = 3
$ws.Cells.Item(159, 21).Value = 1.975
$ws.Cells.Item(159, 22).Value = 1.875
$ws.Cells.Item(159, 23).Value = 0.8
$ws.Cells.Item(159, 25).Value = -1
$ws.Cells.Item(159, 26).Value = 0.825
$ws.Cells.Item(159, 27).Value = -1
$ws.Cells.Item(159, 28).Value = 0
$ws.Cells.Item(159, 29).Value = -0
$ws.Cells.Item(160, 2).Value = 6810145
$ws.Cells.Item(160, 6).Value = 'KV Kortrijk'
$ws.Cells.Item(160, 7).Value = 'Gent'
$ws.Cells.Item(160, 8).Value = 0
$ws.Cells.Item(160, 9).Value = 2
$ws.Cells.Item(160, 10).Value = 'A'
$ws.Cells.Item(160, 11).Value = 7
$ws.Cells.Item(160, 12).Value = 5
$ws.Cells.Item(160, 13).Value = 1.4
$ws.Cells.Item(160, 14).Value = 8.5
$ws.Cells.Item(160, 15).Value = 5.75
$ws.Cells.Item(160, 16).Value = 1.285
$ws.Cells.Item(160, 17).Value = 1.5
$ws.Cells.Item(160, 18).Value = 2.025
$ws.Cells.Item(160, 19).Value = 1.825
$ws.Cells.Item(160, 20).Value = 3.25
$ws.Cells.Item(160, 21).Value = 2.05
$ws.Cells.Item(160, 22).Value = 1.8
$ws.Cells.Item(160, 23).Value = -1
$ws.Cells.Item(160, 25).Value = 0.2849999999999999
$ws.Cells.Item(160, 26).Value = -1
$ws.Cells.Item(160, 27).Value = 0.825
$ws.Cells.Item(160, 28).Value = -1
$ws.Cells.Item(160, 29).Value = 0.8
$ws.Cells.Item(175, 2).Value = 6810159
$ws.Cells.Item(175, 6).Value = 'Charleroi'
$ws.Cells.Item(175, 7).Value = 'Club Brugge'
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(175, 9).Value = 4
$ws.Cells.Item(175, 10).Value = 'A'
$ws.Cells.Item(175, 11).Value = 5.5
$ws.Cells.Item(175, 12).Value = 4.2
$ws.Cells.Item(175, 13).Value = 1.5
$ws.Cells.Item(175, 14).Value = 6
$ws.Cells.Item(175, 15).Value = 4.2
$ws.Cells.Item(175, 16).Value = 1.45
$ws.Cells.Item(175, 17).Value = 1.25
$ws.Cells.Item(175, 18).Value = 1.825
$ws.Cells.Item(175, 19).Value = 2.025
$ws.Cells.Item(175, 20).Value = 2.75
$ws.Cells.Item(175, 21).Value = 1.95
$ws.Cells.Item(175, 22).Value = 1.9
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = 0.45
$ws.Cells.Item(175, 27).Value = 1.025
$ws.Cells.Item(175, 28).Value = 0.95
$ws.Cells.Item(176, 2).Value = 6810158
$ws.Cells.Item(176, 6).Value = 'Gent'
$ws.Cells.Item(176, 7).Value = 'Westerlo'
$ws.Cells.Item(176, 8).Value = 2
$ws.Cells.Item(176, 9).Value = 2
$ws.Cells.Item(176, 10).Value = 'D'
$ws.Cells.Item(176, 11).Value = 1.363
$ws.Cells.Item(176, 12).Value = 5
$ws.Cells.Item(176, 13).Value = 6.5
$ws.Cells.Item(176, 14).Value = 1.444
$ws.Cells.Item(176, 15).Value = 4.5
$ws.Cells.Item(176, 16).Value = 5.75
$ws.Cells.Item(176, 17).Value = -1.25
$ws.Cells.Item(176, 18).Value = 2
$ws.Cells.Item(176, 19).Value = 1.85
$ws.Cells.Item(176, 20).Value = 3
$ws.Cells.Item(176, 21).Value = 2.025
$ws.Cells.Item(176, 22).Value = 1.825
$ws.Cells.Item(176, 24).Value = 3.5
$ws.Cells.Item(176, 25).Value = -1
$ws.Cells.Item(176, 27).Value = 0.8500000000000001
$ws.Cells.Item(176, 28).Value = 1.025
$ws.Cells.Item(183, 2).Value = 6810165
$ws.Cells.Item(183, 6).Value = 'Charleroi'
$ws.Cells.Item(183, 7).Value = 'Eupen'
$ws.Cells.Item(183, 8).Value = 1
$ws.Cells.Item(183, 9).Value = 0
$ws.Cells.Item(183, 11).Value = 1.6
$ws.Cells.Item(183, 12).Value = 4
$ws.Cells.Item(183, 13).Value = 5
$ws.Cells.Item(183, 14).Value = 1.8
$ws.Cells.Item(183, 15).Value = 3.75
$ws.Cells.Item(183, 16).Value = 4.2
$ws.Cells.Item(183, 17).Value = -0.75
$ws.Cells.Item(183, 18).Value = 2.05
$ws.Cells.Item(183, 19).Value = 1.8
$ws.Cells.Item(183, 20).Value = 2.75
$ws.Cells.Item(183, 21).Value = 1.95
$ws.Cells.Item(183, 22).Value = 1.9
$ws.Cells.Item(183, 23).Value = 0.8
$ws.Cells.Item(183, 26).Value = 0.5249999999999999
$ws.Cells.Item(183, 27).Value = -0.5
$ws.Cells.Item(183, 28).Value = -1
$ws.Cells.Item(183, 29).Value = 0.8999999999999999
$ws.Cells.Item(184, 2).Value = 6810168
$ws.Cells.Item(184, 6).Value = 'OH Leuven'
$ws.Cells.Item(184, 7).Value = 'Genk'
$ws.Cells.Item(184, 8).Value = 2
$ws.Cells.Item(184, 9).Value = 1
$ws.Cells.Item(184, 11).Value = 4.5
$ws.Cells.Item(184, 12).Value = 4.2
$ws.Cells.Item(184, 13).Value = 1.666
$ws.Cells.Item(184, 14).Value = 4.333
$ws.Cells.Item(184, 15).Value = 4
$ws.Cells.Item(184, 16).Value = 1.7
$ws.Cells.Item(184, 17).Value = 0.75
$ws.Cells.Item(184, 18).Value = 1.95
$ws.Cells.Item(184, 19).Value = 1.9
$ws.Cells.Item(184, 20).Value = 3
$ws.Cells.Item(184, 21).Value = 1.975
$ws.Cells.Item(184, 22).Value = 1.875
$ws.Cells.Item(184, 23).Value = 3.333
$ws.Cells.Item(184, 26).Value = 0.95
$ws.Cells.Item(184, 27).Value = -1
$ws.Cells.Item(184, 28).Value = 0
$ws.Cells.Item(184, 29).Value = -0
$ws.Cells.Item(185, 2).Value = 6810162
$ws.Cells.Item(185, 6).Value = 'Standard Liege'
$ws.Cells.Item(185, 7).Value = 'Antwerp'
$ws.Cells.Item(185, 8).Value = 0
$ws.Cells.Item(185, 9).Value = 1
$ws.Cells.Item(185, 10).Value = 'A'
$ws.Cells.Item(185, 11).Value = 4
$ws.Cells.Item(185, 12).Value = 3.6
$ws.Cells.Item(185, 13).Value = 1.85
$ws.Cells.Item(185, 14).Value = 3.1
$ws.Cells.Item(185, 15).Value = 3.2
$ws.Cells.Item(185, 16).Value = 2.3
$ws.Cells.Item(185, 17).Value = 0.25
$ws.Cells.Item(185, 18).Value = 1.8
$ws.Cells.Item(185, 19).Value = 2.05
$ws.Cells.Item(185, 20).Value = 2.25
$ws.Cells.Item(185, 21).Value = 1.875
$ws.Cells.Item(185, 22).Value = 1.975
$ws.Cells.Item(185, 23).Value = -1
$ws.Cells.Item(185, 25).Value = 1.3
$ws.Cells.Item(185, 27).Value = 1.05
$ws.Cells.Item(185, 28).Value = -1
$ws.Cells.Item(185, 29).Value = 0.9750000000000001
$ws.Cells.Item(186, 2).Value = 6810164
$ws.Cells.Item(186, 6).Value = 'Union Saint Gilloise'
$ws.Cells.Item(186, 7).Value = 'RWD Molenbeek'
$ws.Cells.Item(186, 8).Value = 3
$ws.Cells.Item(186, 9).Value = 2
$ws.Cells.Item(186, 10).Value = 'H'
$ws.Cells.Item(186, 11).Value = 1.2
$ws.Cells.Item(186, 12).Value = 7
$ws.Cells.Item(186, 13).Value = 12
$ws.Cells.Item(186, 14).Value = 1.142
$ws.Cells.Item(186, 15).Value = 8.5
$ws.Cells.Item(186, 16).Value = 15
$ws.Cells.Item(186, 17).Value = -2.25
$ws.Cells.Item(186, 18).Value = 1.925
$ws.Cells.Item(186, 19).Value = 1.925
$ws.Cells.Item(186, 20).Value = 3.5
$ws.Cells.Item(186, 21).Value = 2.025
$ws.Cells.Item(186, 22).Value = 1.825
$ws.Cells.Item(186, 23).Value = 0.1419999999999999
$ws.Cells.Item(186, 25).Value = -1
$ws.Cells.Item(186, 27).Value = 0.925
$ws.Cells.Item(186, 28).Value = 1.025
$ws.Cells.Item(186, 29).Value = -1
$ws.Cells.Item(229, 8).Value = 3
$ws.Cells.Item(229, 9).Value = 1
$ws.Cells.Item(229, 10).Value = 'H'
$ws.Cells.Item(229, 14).Value = 2.05
$ws.Cells.Item(229, 15).Value = 3.5
$ws.Cells.Item(229, 16).Value = 3.5
$ws.Cells.Item(229, 17).Value = -0.25
$ws.Cells.Item(229, 18).Value = 1.8
$ws.Cells.Item(229, 19).Value = 2.05
$ws.Cells.Item(229, 21).Value = 1.925
$ws.Cells.Item(229, 22).Value = 1.925
$ws.Cells.Item(229, 23).Value = 1.05
$ws.Cells.Item(229, 24).Value = -1
$ws.Cells.Item(229, 25).Value = -1
$ws.Cells.Item(229, 26).Value = 0.8
$ws.Cells.Item(229, 27).Value = -1
$ws.Cells.Item(229, 28).Value = 0.925
$ws.Cells.Item(229, 29).Value = -1
$ws.Cells.Item(230, 14).Value = 5.5
$ws.Cells.Item(230, 15).Value = 4
$ws.Cells.Item(230, 16).Value = 1.6
$ws.Cells.Item(230, 18).Value = 1.85
$ws.Cells.Item(230, 19).Value = 2
$ws.Cells.Item(230, 20).Value = 2.75
$ws.Cells.Item(230, 21).Value = 1.825
$ws.Cells.Item(230, 22).Value = 2.025
$ws.Cells.Item(231, 14).Value = 3.25
$ws.Cells.Item(231, 15).Value = 3.5
$ws.Cells.Item(231, 16).Value = 2.1
$ws.Cells.Item(231, 17).Value = 0.25
$ws.Cells.Item(231, 18).Value = 2.025
$ws.Cells.Item(231, 19).Value = 1.825
$ws.Cells.Item(231, 21).Value = 1.875
$ws.Cells.Item(231, 22).Value = 1.975
$ws.Cells.Item(232, 14).Value = 5.25
$ws.Cells.Item(232, 16).Value = 1.571
$ws.Cells.Item(232, 18).Value = 1.9
$ws.Cells.Item(232, 19).Value = 1.95
$ws.Cells.Item(232, 21).Value = 1.825
$ws.Cells.Item(232, 22).Value = 2.025
$ws.Cells.Item(233, 17).Value = -0.75
$ws.Cells.Item(233, 18).Value = 1.825
$ws.Cells.Item(233, 19).Value = 2.025
$ws.Cells.Item(233, 21).Value = 1.95
$ws.Cells.Item(233, 22).Value = 1.9
$ws.Cells.Item(235, 14).Value = 1.55
$ws.Cells.Item(235, 16).Value = 5.5
$ws.Cells.Item(235, 17).Value = -1
$ws.Cells.Item(235, 18).Value = 2
$ws.Cells.Item(235, 19).Value = 1.85
$ws.Cells.Item(236, 14).Value = 3.3
$ws.Cells.Item(236, 16).Value = 2.15
$ws.Cells.Item(236, 18).Value = 1.975
$ws.Cells.Item(236, 19).Value = 1.875
